$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: convert the text timestamps into real Excel date/time
# serial values, formatted as date-time. Applying a throwaway numeric
# date format first and then the real one mirrors how the source
# workbook ended up with two numFmt entries (164 unused, 165 in use).
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A2").Value = 45687.51764583333
$ws.Range("A3").Value = 45687.52965625
$ws.Range("A4").Value = 45687.52990740741
$ws.Range("A5").Value = 45687.51764467593
$ws.Range("A6").Value = 45687.52965509259
$ws.Range("A7").Value = 45687.52990509259

# --- New column F: "Trening" header (same style as the other headers)
# plus "Gra" for every data row.
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F2").Value = "Gra"
$ws.Range("F3").Value = "Gra"
$ws.Range("F4").Value = "Gra"
$ws.Range("F5").Value = "Gra"
$ws.Range("F6").Value = "Gra"
$ws.Range("F7").Value = "Gra"
